$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Mfge8"
$ws.Cells.Item(2, 3).Value = "Itgb3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 15.12165933333333
$ws.Cells.Item(2, 8).Value = 45.364978
$ws.Cells.Item(2, 9).Value = 0.1696222886509932
$ws.Cells.Item(2, 10).Value = 0.1696222886509932
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.481489333333333
$ws.Cells.Item(2, 14).Value = 7.444467999999999
$ws.Cells.Item(2, 15).Value = 0.2345069082418988
$ws.Cells.Item(2, 16).Value = 0.2345069082418987
$ws.Cells.Item(2, 17).Value = 37.5242363379671
$ws.Cells.Item(2, 18).Value = 337.718127041704
$ws.Cells.Item(2, 19).Value = 0.03977759848045934
$ws.Cells.Item(2, 20).Value = 0.03977759848045932

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Mfge8"
$ws.Cells.Item(3, 3).Value = "Itgb3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 15.12165933333333
$ws.Cells.Item(3, 8).Value = 45.364978
$ws.Cells.Item(3, 9).Value = 0.1696222886509932
$ws.Cells.Item(3, 10).Value = 0.1696222886509932
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.245227
$ws.Cells.Item(3, 14).Value = 21.735681
$ws.Cells.Item(3, 15).Value = 0.6846919551326144
$ws.Cells.Item(3, 16).Value = 0.6846919551326142
$ws.Cells.Item(3, 17).Value = 109.5598544866687
$ws.Cells.Item(3, 18).Value = 986.038690380018
$ws.Cells.Item(3, 19).Value = 0.1161390164505172
$ws.Cells.Item(3, 20).Value = 0.1161390164505172

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Mfge8"
$ws.Cells.Item(4, 3).Value = "Itgb3"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 15.12165933333333
$ws.Cells.Item(4, 8).Value = 45.364978
$ws.Cells.Item(4, 9).Value = 0.1696222886509932
$ws.Cells.Item(4, 10).Value = 0.1696222886509932
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.2001876666666667
$ws.Cells.Item(4, 14).Value = 0.600563
$ws.Cells.Item(4, 15).Value = 0.01891823194544989
$ws.Cells.Item(4, 16).Value = 0.01891823194544989
$ws.Cells.Item(4, 17).Value = 3.027169698068222
$ws.Cells.Item(4, 18).Value = 27.244527282614
$ws.Cells.Item(4, 19).Value = 0.003208953799817543
$ws.Cells.Item(4, 20).Value = 0.003208953799817542

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Mfge8"
$ws.Cells.Item(5, 3).Value = "Itgb3"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 15.12165933333333
$ws.Cells.Item(5, 8).Value = 45.364978
$ws.Cells.Item(5, 9).Value = 0.1696222886509932
$ws.Cells.Item(5, 10).Value = 0.1696222886509932
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6548283333333333
$ws.Cells.Item(5, 14).Value = 1.964485
$ws.Cells.Item(5, 15).Value = 0.06188290468003712
$ws.Cells.Item(5, 16).Value = 0.06188290468003711
$ws.Cells.Item(5, 17).Value = 9.902090978481111
$ws.Cells.Item(5, 18).Value = 89.11881880633
$ws.Cells.Item(5, 19).Value = 0.01049671992019916
$ws.Cells.Item(5, 20).Value = 0.01049671992019915

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Mfge8"
$ws.Cells.Item(6, 3).Value = "Itgb3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 30.93224
$ws.Cells.Item(6, 8).Value = 92.79671999999999
$ws.Cells.Item(6, 9).Value = 0.3469723279862584
$ws.Cells.Item(6, 10).Value = 0.3469723279862583
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.481489333333333
$ws.Cells.Item(6, 14).Value = 7.444467999999999
$ws.Cells.Item(6, 15).Value = 0.2345069082418988
$ws.Cells.Item(6, 16).Value = 0.2345069082418987
$ws.Cells.Item(6, 17).Value = 76.75802361610664
$ws.Cells.Item(6, 18).Value = 690.8222125449598
$ws.Cells.Item(6, 19).Value = 0.0813674078815515
$ws.Cells.Item(6, 20).Value = 0.08136740788155146

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Mfge8"
$ws.Cells.Item(7, 3).Value = "Itgb3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 30.93224
$ws.Cells.Item(7, 8).Value = 92.79671999999999
$ws.Cells.Item(7, 9).Value = 0.3469723279862584
$ws.Cells.Item(7, 10).Value = 0.3469723279862583
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.245227
$ws.Cells.Item(7, 14).Value = 21.735681
$ws.Cells.Item(7, 15).Value = 0.6846919551326144
$ws.Cells.Item(7, 16).Value = 0.6846919551326142
$ws.Cells.Item(7, 17).Value = 224.11110041848
$ws.Cells.Item(7, 18).Value = 2016.99990376632
$ws.Cells.Item(7, 19).Value = 0.237569161625826
$ws.Cells.Item(7, 20).Value = 0.2375691616258259

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Mfge8"
$ws.Cells.Item(8, 3).Value = "Itgb3"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 30.93224
$ws.Cells.Item(8, 8).Value = 92.79671999999999
$ws.Cells.Item(8, 9).Value = 0.3469723279862584
$ws.Cells.Item(8, 10).Value = 0.3469723279862583
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.2001876666666667
$ws.Cells.Item(8, 14).Value = 0.600563
$ws.Cells.Item(8, 15).Value = 0.01891823194544989
$ws.Cells.Item(8, 16).Value = 0.01891823194544989
$ws.Cells.Item(8, 17).Value = 6.192252950373332
$ws.Cells.Item(8, 18).Value = 55.73027655335999
$ws.Cells.Item(8, 19).Value = 0.006564102979496751
$ws.Cells.Item(8, 20).Value = 0.006564102979496749

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Mfge8"
$ws.Cells.Item(9, 3).Value = "Itgb3"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 30.93224
$ws.Cells.Item(9, 8).Value = 92.79671999999999
$ws.Cells.Item(9, 9).Value = 0.3469723279862584
$ws.Cells.Item(9, 10).Value = 0.3469723279862583
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.6548283333333333
$ws.Cells.Item(9, 14).Value = 1.964485
$ws.Cells.Item(9, 15).Value = 0.06188290468003712
$ws.Cells.Item(9, 16).Value = 0.06188290468003711
$ws.Cells.Item(9, 17).Value = 20.25530716546666
$ws.Cells.Item(9, 18).Value = 182.2977644892
$ws.Cells.Item(9, 19).Value = 0.0214716554993842
$ws.Cells.Item(9, 20).Value = 0.0214716554993842

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Mfge8"
$ws.Cells.Item(10, 3).Value = "Itgb3"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.447555666666666
$ws.Cells.Item(10, 8).Value = 13.342667
$ws.Cells.Item(10, 9).Value = 0.04988900718188559
$ws.Cells.Item(10, 10).Value = 0.04988900718188558
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.481489333333333
$ws.Cells.Item(10, 14).Value = 7.444467999999999
$ws.Cells.Item(10, 15).Value = 0.2345069082418988
$ws.Cells.Item(10, 16).Value = 0.2345069082418987
$ws.Cells.Item(10, 17).Value = 11.03656194623955
$ws.Cells.Item(10, 18).Value = 99.32905751615597
$ws.Cells.Item(10, 19).Value = 0.01169931682948187
$ws.Cells.Item(10, 20).Value = 0.01169931682948187

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Mfge8"
$ws.Cells.Item(11, 3).Value = "Itgb3"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.447555666666666
$ws.Cells.Item(11, 8).Value = 13.342667
$ws.Cells.Item(11, 9).Value = 0.04988900718188559
$ws.Cells.Item(11, 10).Value = 0.04988900718188558
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 7.245227
$ws.Cells.Item(11, 14).Value = 21.735681
$ws.Cells.Item(11, 15).Value = 0.6846919551326144
$ws.Cells.Item(11, 16).Value = 0.6846919551326142
$ws.Cells.Item(11, 17).Value = 32.22355040013633
$ws.Cells.Item(11, 18).Value = 290.011953601227
$ws.Cells.Item(11, 19).Value = 0.03415860186699028
$ws.Cells.Item(11, 20).Value = 0.03415860186699027

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Mfge8"
$ws.Cells.Item(12, 3).Value = "Itgb3"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.447555666666666
$ws.Cells.Item(12, 8).Value = 13.342667
$ws.Cells.Item(12, 9).Value = 0.04988900718188559
$ws.Cells.Item(12, 10).Value = 0.04988900718188558
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.2001876666666667
$ws.Cells.Item(12, 14).Value = 0.600563
$ws.Cells.Item(12, 15).Value = 0.01891823194544989
$ws.Cells.Item(12, 16).Value = 0.01891823194544989
$ws.Cells.Item(12, 17).Value = 0.8903457912801109
$ws.Cells.Item(12, 18).Value = 8.013112121520999
$ws.Cells.Item(12, 19).Value = 0.0009438118093951271
$ws.Cells.Item(12, 20).Value = 0.0009438118093951268

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Mfge8"
$ws.Cells.Item(13, 3).Value = "Itgb3"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.447555666666666
$ws.Cells.Item(13, 8).Value = 13.342667
$ws.Cells.Item(13, 9).Value = 0.04988900718188559
$ws.Cells.Item(13, 10).Value = 0.04988900718188558
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.6548283333333333
$ws.Cells.Item(13, 14).Value = 1.964485
$ws.Cells.Item(13, 15).Value = 0.06188290468003712
$ws.Cells.Item(13, 16).Value = 0.06188290468003711
$ws.Cells.Item(13, 17).Value = 2.912385464610555
$ws.Cells.Item(13, 18).Value = 26.211469181495
$ws.Cells.Item(13, 19).Value = 0.003087276676018313
$ws.Cells.Item(13, 20).Value = 0.003087276676018312

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Mfge8"
$ws.Cells.Item(14, 3).Value = "Itgb3"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 38.64755633333333
$ws.Cells.Item(14, 8).Value = 115.942669
$ws.Cells.Item(14, 9).Value = 0.4335163761808628
$ws.Cells.Item(14, 10).Value = 0.4335163761808628
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 2.481489333333333
$ws.Cells.Item(14, 14).Value = 7.444467999999999
$ws.Cells.Item(14, 15).Value = 0.2345069082418988
$ws.Cells.Item(14, 16).Value = 0.2345069082418987
$ws.Cells.Item(14, 17).Value = 95.90349880056576
$ws.Cells.Item(14, 18).Value = 863.131489205092
$ws.Cells.Item(14, 19).Value = 0.1016625850504061
$ws.Cells.Item(14, 20).Value = 0.101662585050406

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Mfge8"
$ws.Cells.Item(15, 3).Value = "Itgb3"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 38.64755633333333
$ws.Cells.Item(15, 8).Value = 115.942669
$ws.Cells.Item(15, 9).Value = 0.4335163761808628
$ws.Cells.Item(15, 10).Value = 0.4335163761808628
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 7.245227
$ws.Cells.Item(15, 14).Value = 21.735681
$ws.Cells.Item(15, 15).Value = 0.6846919551326144
$ws.Cells.Item(15, 16).Value = 0.6846919551326142
$ws.Cells.Item(15, 17).Value = 280.0103186302877
$ws.Cells.Item(15, 18).Value = 2520.092867672589
$ws.Cells.Item(15, 19).Value = 0.2968251751892809
$ws.Cells.Item(15, 20).Value = 0.2968251751892808

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Mfge8"
$ws.Cells.Item(16, 3).Value = "Itgb3"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 38.64755633333333
$ws.Cells.Item(16, 8).Value = 115.942669
$ws.Cells.Item(16, 9).Value = 0.4335163761808628
$ws.Cells.Item(16, 10).Value = 0.4335163761808628
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.2001876666666667
$ws.Cells.Item(16, 14).Value = 0.600563
$ws.Cells.Item(16, 15).Value = 0.01891823194544989
$ws.Cells.Item(16, 16).Value = 0.01891823194544989
$ws.Cells.Item(16, 17).Value = 7.736764124738555
$ws.Cells.Item(16, 18).Value = 69.63087712264701
$ws.Cells.Item(16, 19).Value = 0.008201363356740472
$ws.Cells.Item(16, 20).Value = 0.008201363356740468

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Mfge8"
$ws.Cells.Item(17, 3).Value = "Itgb3"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 38.64755633333333
$ws.Cells.Item(17, 8).Value = 115.942669
$ws.Cells.Item(17, 9).Value = 0.4335163761808628
$ws.Cells.Item(17, 10).Value = 0.4335163761808628
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.6548283333333333
$ws.Cells.Item(17, 14).Value = 1.964485
$ws.Cells.Item(17, 15).Value = 0.06188290468003712
$ws.Cells.Item(17, 16).Value = 0.06188290468003711
$ws.Cells.Item(17, 17).Value = 25.30751490116278
$ws.Cells.Item(17, 18).Value = 227.767634110465
$ws.Cells.Item(17, 19).Value = 0.02682725258443545
$ws.Cells.Item(17, 20).Value = 0.02682725258443544
